$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票 (stock) sheet

# 1) Fix the typo in the company-name shared string (drop the stray space)
$ws.Range("B2").Value = "力宇創業投資股份有限公司"

# 2) Add a new "property_category" column header at the end of row 1 (K1),
#    matching the bold/bordered header style used by the rest of row 1.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "property_category"

# 3) Record the corresponding "stock" value for the data row, inserted at I2
#    (ahead of the legislator_name / legislator_id values, which shift right
#    by one column: I2->J2, J2->K2).
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("K2").Value = $ws.Range("J2").Value()
$ws.Range("J2").Value = $ws.Range("I2").Value()
$ws.Range("I2").Value = "stock"
